$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.090710333333333
$ws.Range("N2").Value = 3.272131
$ws.Range("O2").Value = 0.0488470045579656
$ws.Range("P2").Value = 0.0488470045579656
$ws.Range("Q2").Value = 0.2296221564951111
$ws.Range("R2").Value = 2.066599408456
$ws.Range("S2").Value = 0.0488470045579656
$ws.Range("T2").Value = 0.0488470045579656

# Row 3
$ws.Range("O3").Value = 0.7616320856558244
$ws.Range("P3").Value = 0.7616320856558244
$ws.Range("S3").Value = 0.7616320856558244
$ws.Range("T3").Value = 0.7616320856558244

# Row 4
$ws.Range("M4").Value = 0.740281
$ws.Range("N4").Value = 2.220843
$ws.Range("O4").Value = 0.03315317392351528
$ws.Range("P4").Value = 0.03315317392351528
$ws.Range("Q4").Value = 0.1558479042853333
$ws.Range("R4").Value = 1.402631138568
$ws.Range("S4").Value = 0.03315317392351528
$ws.Range("T4").Value = 0.03315317392351528

# Row 5
$ws.Range("M5").Value = 2.784013333333333
$ws.Range("N5").Value = 8.352039999999999
$ws.Range("O5").Value = 0.1246808688124989
$ws.Range("P5").Value = 0.1246808688124989
$ws.Range("Q5").Value = 0.5861053350044444
$ws.Range("R5").Value = 5.27494801504
$ws.Range("S5").Value = 0.1246808688124989
$ws.Range("T5").Value = 0.1246808688124989

# Row 6
$ws.Range("M6").Value = 0.2710316666666667
$ws.Range("N6").Value = 0.813095
$ws.Range("O6").Value = 0.01213803945228936
$ws.Range("P6").Value = 0.01213803945228936
$ws.Range("Q6").Value = 0.05705903196888889
$ws.Range("R6").Value = 0.51353128772
$ws.Range("S6").Value = 0.01213803945228936
$ws.Range("T6").Value = 0.01213803945228936

# Row 7
$ws.Range("M7").Value = 0.436508
$ws.Range("N7").Value = 1.309524
$ws.Range("O7").Value = 0.01954882759790648
$ws.Range("P7").Value = 0.01954882759790648
$ws.Range("Q7").Value = 0.09189599220266666
$ws.Range("R7").Value = 0.827063929824
$ws.Range("S7").Value = 0.01954882759790648
$ws.Range("T7").Value = 0.01954882759790648
